# BOT; UPDATE DATA (#1897)
# Appends two new daily rows (2020-05-24 / serial 43975, and
# 2020-05-25 / serial 43976) to the "相談件数" sheet. The existing
# "※..." total/footnote row shifts from row 120 down to row 122.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("相談件数")

# Shift the footer row (currently row 120) down two rows to make room
# for the two new data rows; Insert() carries the row-120 formatting
# down onto rows 120/121 (matching the existing data rows' styles).
$ws.Range("A120:E121").Insert()

# New row 120: 2020-05-24
$ws.Cells.Item(120, 1).Value = 43975
$ws.Cells.Item(120, 2).Value = 111
$ws.Cells.Item(120, 3).Value = 38846
$ws.Cells.Item(120, 4).Value = 0
$ws.Cells.Item(120, 5).Value = 7801

# New row 121: 2020-05-25
$ws.Cells.Item(121, 1).Value = 43976
$ws.Cells.Item(121, 2).Value = 176
$ws.Cells.Item(121, 3).Value = 39022
$ws.Cells.Item(121, 4).Value = 41
$ws.Cells.Item(121, 5).Value = 7842

# Update the workbook-level print area defined name to cover the new
# bottom row (122).
$wb.Names.Item("相談件数!Print_Area").RefersTo = "=相談件数!`$A`$1:`$E`$122"

# Update the remembered selection to the new last data cell.
$ws.Activate()
$ws.Range("C128").Select()
